# Update "想去人数" (number interested) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 763
$ws1.Range("F3").Value = 22
$ws1.Range("F5").Value = 34
$ws1.Range("F6").Value = 259
$ws1.Range("F7").Value = 3610
$ws1.Range("F9").Value = 4231
$ws1.Range("F11").Value = 1053
$ws1.Range("F12").Value = 53

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 763
$ws4.Range("F3").Value = 22
$ws4.Range("F5").Value = 34
$ws4.Range("F7").Value = 259
$ws4.Range("F8").Value = 3610
$ws4.Range("F10").Value = 4231
$ws4.Range("F12").Value = 1053
$ws4.Range("F13").Value = 53
